# partial update import asset
# - adds a new "Asset Class ID" lookup sheet (between "Asset Model ID" and "Component ID")
# - reworks the "Upload" sheet header row (adds " *" markers, drops "Warranty expired date",
#   renames "IP computer" -> "Computer IP", adds a "sn" column next to every component id column)
# - freezes the header row on "Upload"
# - refreshes the instructions on "Panduan"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# helper: apply the "header" look (bold white text on navy fill, centered
# vertically) that every lookup-sheet / Upload header row uses in this
# workbook.
# ---------------------------------------------------------------------------
function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 11
    $rng.Font.Color = 16777215      # white
    $rng.Interior.Color = 8388608   # navy (BGR 00 00 80 -> RGB 000080)
    $rng.Interior.Pattern = 1       # xlSolid
    $rng.VerticalAlignment = -4108  # xlCenter
}

# helper: size a column from its header text the same way this workbook's
# generator does ( (chars + 3) * 1.65 "Excel characters" ), translated into
# whatever ColumnWidth value makes the engine's pixel round-trip land as
# close as possible to that target.
function Set-BestFitWidth($ws, $colIndex, $text) {
    $target = ([double]$text.Length + 3) * 1.65
    $px = [math]::Round($target * 7 - 5)
    $cw = $px / 7
    $ws.Columns.Item($colIndex).ColumnWidth = $cw
}

# ===========================================================================
# 1. New "Asset Class ID" lookup sheet, inserted right after "Asset Model ID"
#    (i.e. right before "Component ID").
# ===========================================================================
$assetModelSheet = $wb.Worksheets.Item("Asset Model ID")
$classSheet = $wb.Worksheets.Add($null, $assetModelSheet)
$classSheet.Name = "Asset Class ID"

$classRows = @(
    @("Asset class id", "Name"),
    @("NBADVUSR", "Notebook Advance User"),
    @("NBADVUSR2", "Notebook Advance User 2"),
    @("NBPERFUSR", "Notebook Performance User"),
    @("NBREGUSR", "Notebook Reg. User"),
    @("NBREGUSR2", "Notebook Reg. User 2"),
    @("PERFPC", "Performance PC"),
    @("PERFPC1", "Performance PC 1"),
    @("PERFPC2", "Performance PC 2"),
    @("STDPC1", "Standard PC 1")
)

for ($i = 0; $i -lt $classRows.Length; $i++) {
    $r = $i + 1
    $classSheet.Cells.Item($r, 1).Value = $classRows[$i][0]
    $classSheet.Cells.Item($r, 2).Value = $classRows[$i][1]
}
# last data row: "safety net" has a numeric 0 id
$classSheet.Cells.Item(11, 1).Value = 0
$classSheet.Cells.Item(11, 2).Value = "safety net"

Set-HeaderStyle($classSheet.Range("A1:B1"))

Set-BestFitWidth $classSheet 1 "Asset class id"
$classSheet.Columns.Item(2).ColumnWidth = 27.142857142857142   # best-fit on "Notebook Advance User 2"

$classSheet.Range("A1").Select()

# ===========================================================================
# 2. "Upload" sheet: rebuild the header row.
# ===========================================================================
$uploadSheet = $wb.Worksheets.Item("Upload")

$uploadHeaders = @(
    "Tagging id *", "Project id *", "Site id *", "Asset model id *", "Asset class id",
    "DO number", "Computer name", "Computer IP", "CPU sn", "Monitor sn",
    "Keyboard sn", "Shipping date", "Description", "Mouse id", "Mouse sn",
    "Floopy disk id", "Floopy disk sn", "Processor id", "Processor sn", "Memory id",
    "Memory sn", "Hardisk id", "Hardisk sn", "CD / DVD rom id", "CD / DVD rom sn",
    "NIC id", "NIC sn", "Others id", "Others sn"
)

# clear out the previous 22-column header first so no stale cells are left
# to the right once the sheet only needs 29 (a superset, but keep this
# robust if it is ever re-run).
$uploadSheet.Range("A1:AZ1").ClearContents()

for ($i = 0; $i -lt $uploadHeaders.Length; $i++) {
    $col = $i + 1
    $uploadSheet.Cells.Item(1, $col).Value = $uploadHeaders[$i]
    Set-BestFitWidth $uploadSheet $col $uploadHeaders[$i]
}

$headerRange = $uploadSheet.Range($uploadSheet.Cells.Item(1, 1), $uploadSheet.Cells.Item(1, $uploadHeaders.Length))
Set-HeaderStyle($headerRange)

# freeze the header row
$uploadSheet.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ===========================================================================
# 3. "Panduan" sheet: refresh the instructions text + add a 6th line.
# ===========================================================================
$panduanSheet = $wb.Worksheets.Item("Panduan")

$panduanSheet.Range("A1").Value = "Panduan upload"
$panduanSheet.Range("A2").Value = "1. Lengkapi semua data-data yang ada pada sheet Upload"
$panduanSheet.Range("A3").Value = "2. Kolom pada sheet Upload dengan simbol (*) artinya wajib diisi"
$panduanSheet.Range("A4").Value = "3. Tagging id harus unik (tidak boleh sama)"
$panduanSheet.Range("A5").Value = "4. Kolom Project id, Site id, Asset model id, Asset class id diisi dengan id masing-masing. Id bisa dicek pada masing-masing sheet sesuai nama kolom"
$panduanSheet.Range("A6").Value = "5. Kolom Mouse id, Floopy disk id, Processor id, Memory id, Hardisk id, CD / DVD room id, NIC id, Other id diisi dengan id masing-masing. Id bisa dicek pada sheet ``Components ID``"

Set-HeaderStyle($panduanSheet.Range("A1"))
Set-BestFitWidth $panduanSheet 1 $panduanSheet.Range("A6").Value

$uploadSheet.Activate()
$uploadSheet.Range("A1").Select()

Write-Output "asset-import template updated"
